# FIX: Corrected stoichiometric factor in Streeter-Phelps example
#
# The "stoi" (stoichiometry) sheet's "reox" row had its "s" factor
# (column C, row 3) mistakenly pointing at the text "ka" instead of
# the correct numeric stoichiometric factor of 1. Also bring the
# active/selected sheet in the workbook view back to "stoi" (it was
# left on "pars").

$wb = $excel.ActiveWorkbook

$stoiSheet = $wb.Worksheets.Item("stoi")

# Correct the stoichiometric factor for the "reox" row: it referenced
# the text "ka" by mistake; it should be the number 1.
$stoiSheet.Range("C3").Value = 1

# The workbook should re-open showing the "stoi" sheet as active/selected
# (previously "pars" was the selected tab).
$stoiSheet.Activate()

# Match the saved window chrome (tab-bar / horizontal-scrollbar split ratio).
$win = $excel.ActiveWindow
$win.TabRatio = 0.157
